# Add "cluster_class" column (uncertain DBSCAN clustering method) between
# "work" and "shannon_entropy" on the aucs degree-centrality results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "shannon_entropy" column (G) one position to the right
# so the new column can be inserted in its place.
$ws.Columns("G:G").Insert()

# Header for the newly inserted column.
$ws.Range("G1").Value = "cluster_class"

# cluster_class values for each data row (A2:A63 -> U-labels, A63 is "mean").
$values = @(0,0,0,0,0,0,0,0,0,1,0,0,1,1,0,1,1,1,1,0,1,-1,1,0,1,0,1,1,1,0,0,0,1,-1,0,1,0,0,1,0,1,0,0,0,1,1,0,0,1,0,0,1,1,0,0,1,1,-1,1,1,1,0.33)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
